$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Formula = '''2026-02-23 21:18:34'
$ws.Range('O2').Formula = '''5.8 °C'
$ws.Range('E3').Formula = '''2026-02-23 21:18:36'
$ws.Range('K3').Formula = '''16.6 MJ/m2'
$ws.Range('E4').Formula = '''2026-02-23 21:18:38'
$ws.Range('H4').Formula = '''69%'
$ws.Range('O4').Formula = '''12.1 °C'
$ws.Range('E5').Formula = '''2026-02-23 21:18:41'
$ws.Range('H5').Formula = '''29%'
$ws.Range('E6').Formula = '''2026-02-23 21:18:43'
$ws.Range('E7').Formula = '''2026-02-23 21:18:45'
$ws.Range('H7').Formula = '''67%'
$ws.Range('E8').Formula = '''2026-02-23 21:18:48'
$ws.Range('E9').Formula = '''2026-02-23 21:18:50'
$ws.Range('H9').Formula = '''73%'
$ws.Range('O9').Formula = '''12.5 °C'
$ws.Range('E10').Formula = '''2026-02-23 21:18:53'
$ws.Range('K10').Formula = '''15.3 MJ/m2'
$ws.Range('O10').Formula = '''10.8 °C'
$ws.Range('E11').Formula = '''2026-02-23 21:18:55'
$ws.Range('E12').Formula = '''2026-02-23 21:18:57'
$ws.Range('H12').Formula = '''85%'
$ws.Range('E13').Formula = '''2026-02-23 21:18:59'
$ws.Range('H13').Formula = '''59%'
$ws.Range('O13').Formula = '''7.1 °C'
$ws.Range('E14').Formula = '''2026-02-23 21:19:02'
$ws.Range('N14').Formula = '''6.4 °C 20:35 TU'
$ws.Range('O14').Formula = '''12.5 °C'
$ws.Range('E15').Formula = '''2026-02-23 21:19:04'
$ws.Range('O15').Formula = '''12.5 °C'
$ws.Range('E16').Formula = '''2026-02-23 21:19:06'
$ws.Range('H16').Formula = '''22%'
$ws.Range('O16').Formula = '''3.8 °C'
$ws.Range('E17').Formula = '''2026-02-23 21:19:09'
$ws.Range('E18').Formula = '''2026-02-23 21:19:11'
$ws.Range('O18').Formula = '''11.0 °C'
$ws.Range('E19').Formula = '''2026-02-23 21:19:14'
$ws.Range('E20').Formula = '''2026-02-23 21:19:16'
$ws.Range('E21').Formula = '''2026-02-23 21:19:19'
$ws.Range('E22').Formula = '''2026-02-23 21:19:21'
$ws.Range('H22').Formula = '''24%'
$ws.Range('E23').Formula = '''2026-02-23 21:19:23'
$ws.Range('E24').Formula = '''2026-02-23 21:19:26'
$ws.Range('J24').Formula = '''1026.0 hPa'
$ws.Range('K24').Formula = '''16.1 MJ/m2'
$ws.Range('O24').Formula = '''8.6 °C'
$ws.Range('E25').Formula = '''2026-02-23 21:19:28'
$ws.Range('E26').Formula = '''2026-02-23 21:19:31'
$ws.Range('O26').Formula = '''10.0 °C'
$ws.Range('E27').Formula = '''2026-02-23 21:19:33'
$ws.Range('H27').Formula = '''29%'
$ws.Range('O27').Formula = '''5.7 °C'
$ws.Range('E28').Formula = '''2026-02-23 21:19:36'
$ws.Range('O28').Formula = '''11.1 °C'
$ws.Range('E29').Formula = '''2026-02-23 21:19:38'
$ws.Range('H29').Formula = '''83%'
$ws.Range('O29').Formula = '''10.8 °C'
$ws.Range('E30').Formula = '''2026-02-23 21:19:40'
$ws.Range('E31').Formula = '''2026-02-23 21:19:43'
$ws.Range('O31').Formula = '''16.6 °C'
$ws.Range('E32').Formula = '''2026-02-23 21:19:45'
$ws.Range('H32').Formula = '''67%'
$ws.Range('O32').Formula = '''7.5 °C'
$ws.Range('E33').Formula = '''2026-02-23 21:19:48'
$ws.Range('J33').Formula = '''1025.2 hPa'
$ws.Range('E34').Formula = '''2026-02-23 21:19:50'
$ws.Range('E35').Formula = '''2026-02-23 21:19:53'
$ws.Range('H35').Formula = '''38%'
$ws.Range('O35').Formula = '''12.2 °C'
$ws.Range('E36').Formula = '''2026-02-23 21:19:55'
$ws.Range('E37').Formula = '''2026-02-23 21:19:57'
$ws.Range('J37').Formula = '''1026.7 hPa'
$ws.Range('O37').Formula = '''9.1 °C'
$ws.Range('E38').Formula = '''2026-02-23 21:20:00'
$ws.Range('E39').Formula = '''2026-02-23 21:20:02'
$ws.Range('E40').Formula = '''2026-02-23 21:20:04'
$ws.Range('J40').Formula = '''1026.2 hPa'
$ws.Range('O40').Formula = '''8.8 °C'
$ws.Range('E41').Formula = '''2026-02-23 21:20:06'
$ws.Range('H41').Formula = '''73%'
$ws.Range('O41').Formula = '''12.0 °C'
$ws.Range('E42').Formula = '''2026-02-23 21:20:09'
$ws.Range('O42').Formula = '''11.8 °C'
$ws.Range('E43').Formula = '''2026-02-23 21:20:11'
$ws.Range('E44').Formula = '''2026-02-23 21:20:14'
$ws.Range('O44').Formula = '''3.1 °C'
$ws.Range('E45').Formula = '''2026-02-23 21:20:16'
$ws.Range('J45').Formula = '''1026.9 hPa'
$ws.Range('O45').Formula = '''8.3 °C'
$ws.Range('E46').Formula = '''2026-02-23 21:20:18'
$ws.Range('H46').Formula = '''73%'
